$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.631.30'
$ws.Range("E2").Value = '  +4.18%  '
$ws.Range("D3").Value = '1.747.96'
$ws.Range("E3").Value = '  +4.68%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.90'
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4813'
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  +2.83%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06264'
$ws.Range("E9").Value = '  +1.47%  '
$ws.Range("D10").Value = '1.743.82'
$ws.Range("E10").Value = '  +4.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07113'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.82'
$ws.Range("E12").Value = '  +6.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6177'
$ws.Range("E13").Value = '  +4.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.514'
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.30'
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '26.626.10'
$ws.Range("E17").Value = '  +4.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006912'
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.73'
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("D21").Value = '1.970.73'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.654'
$ws.Range("E22").Value = '  +4.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.864'
$ws.Range("E23").Value = '  +1.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.359'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.20'
$ws.Range("E25").Value = '  -0.31%  '
$ws.Range("E26").Value = '  +2.91%  '
$ws.Range("E27").Value = '  +5.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.410'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '107.85'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.028'
$ws.Range("E30").Value = '  +1.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.779'
$ws.Range("E31").Value = '  +3.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07901'
$ws.Range("E32").Value = '  +0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04586'
$ws.Range("E33").Value = '  +8.67%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9993'
$ws.Range("E35").Value = '  +4.55%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6356'
$ws.Range("E36").Value = '  +4.23%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9507'
$ws.Range("E37").Value = '  +10.38%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '114.49'
$ws.Range("E38").Value = '  +18.48%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.466'
$ws.Range("E39").Value = '  -4.85%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.986'
$ws.Range("E40").Value = '  +6.20%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.004'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01514'
$ws.Range("E42").Value = '  +2.35%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.696'
$ws.Range("E43").Value = '  +16.69%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3921'
$ws.Range("E44").Value = '  +4.02%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.752'
$ws.Range("E45").Value = '  +8.56%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1204'
$ws.Range("E46").Value = '  +7.84%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05329'
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.989'
$ws.Range("E48").Value = '  +7.92%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.93'
$ws.Range("E49").Value = '  +3.40%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.257'
$ws.Range("E50").Value = '  +4.57%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3459'
$ws.Range("E51").Value = '  +3.69%  '

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
